# Apply updates described by the diff across three worksheets:
#   展览 (sheet1), 演出 (sheet2), 全部类型 (sheet4)

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 387
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 678
$ws1.Range("G3").Value = 29
$ws1.Range("F5").Value = 2111
$ws1.Range("F7").Value = 10983
$ws1.Range("F9").Value = 166
$ws1.Range("F12").Value = 10828
$ws1.Range("F13").Value = 438
$ws1.Range("F14").Value = 1130
$ws1.Range("F17").Value = 5427
$ws1.Range("F18").Value = 81
$ws1.Range("F19").Value = 3403

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 387
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 678
$ws4.Range("G3").Value = 29
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F6").Value = 2111
$ws4.Range("F10").Value = 10983
$ws4.Range("F12").Value = 166
$ws4.Range("F15").Value = 10828
$ws4.Range("F16").Value = 438
$ws4.Range("F17").Value = 1130
$ws4.Range("F20").Value = 5427
$ws4.Range("F21").Value = 81
$ws4.Range("F22").Value = 3403
